# Insert two new weekly-report rows (133:134) above the existing data.
# Excel shifts the previously-existing rows 133:183 down to 135:185,
# and we populate the two freshly inserted rows with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("133:134").Insert()

# Row 133 - "Primera" quality, new date 2021-09-27 (serial 44466)
$ws.Cells.Item(133, 1).Value = 1
$ws.Cells.Item(133, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(133, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(133, 4).Value = 44466
$ws.Cells.Item(133, 5).Value = 15
$ws.Cells.Item(133, 6).Value = 100112043
$ws.Cells.Item(133, 7).Value = "Pepino ensalada"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 120
$ws.Cells.Item(133, 11).Value = 10000
$ws.Cells.Item(133, 12).Value = 11000
$ws.Cells.Item(133, 13).Value = 10500
$ws.Cells.Item(133, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(133, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(133, 16).Value = 150
$ws.Cells.Item(133, 17).Value = 70
$ws.Cells.Item(133, 18).Value = "Hortaliza"

# Row 134 - "Segunda" quality, new date 2021-09-27 (serial 44466)
$ws.Cells.Item(134, 1).Value = 1
$ws.Cells.Item(134, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(134, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(134, 4).Value = 44466
$ws.Cells.Item(134, 5).Value = 15
$ws.Cells.Item(134, 6).Value = 100112043
$ws.Cells.Item(134, 7).Value = "Pepino ensalada"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Segunda"
$ws.Cells.Item(134, 10).Value = 160
$ws.Cells.Item(134, 11).Value = 7000
$ws.Cells.Item(134, 12).Value = 8000
$ws.Cells.Item(134, 13).Value = 7500
$ws.Cells.Item(134, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(134, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(134, 16).Value = 75
$ws.Cells.Item(134, 17).Value = 100
$ws.Cells.Item(134, 18).Value = "Hortaliza"
